# Canada Premier League workbook update
# The source data feed re-sorted same-day fixtures, which swapped the row
# order for several pairs of matches that share a round/date grouping.
# Net effect: the entire row content (id, date, teams, score, result,
# odds, P&L columns -- i.e. columns B through AB) for each of the pairs
# below is exchanged between the two rows, while the leading index column
# (A) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(3, 4),
    @(20, 21),
    @(26, 27),
    @(31, 32),
    @(36, 37),
    @(47, 48),
    @(55, 56),
    @(58, 59),
    @(62, 63),
    @(64, 65),
    @(81, 82),
    @(83, 84),
    @(86, 87),
    @(98, 99)
)

for ($i = 0; $i -lt $rowPairs.Count; $i++) {
    $pair = $rowPairs[$i]
    $r1 = $pair[0]
    $r2 = $pair[1]

    $addr1 = "B${r1}:AB${r1}"
    $addr2 = "B${r2}:AB${r2}"

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
